$d = $word.ActiveDocument

# 1) Fix "objectivo" -> "objetivo" (only the first occurrence, the one right
#    after "não foram estabelecidas. O principal "). The search string spans
#    across the run that precedes it so that, just like real Word
#    find/replace across multiple runs, the touched runs collapse into a
#    single run and the now-redundant w:proofErr spell-check markers around
#    "objectivo"/"objetivo" disappear (that word is correctly spelled in
#    pt-BR, unlike the old "objectivo" Portugal-Portuguese spelling).
$old1 = " não foram estabelecidas. O principal objectivo deste"
$new1 = " não foram estabelecidas. O principal objetivo deste"

$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 1)

# 2) Fix the broken hyphenation "prev- alente" -> "prevalente". The search
#    string spans the whole original run group so that, just like real
#    Word find/replace across multiple runs, the matched runs collapse into
#    a single run (and the now-redundant w:proofErr spell-check markers for
#    "prev" disappear).
$old = "O delírio é um estado de confusão aguda caracterizado por declínios na atenção, na consciência e na cognição. As flutuações do estado mental ao longo do tempo são características e necessárias para o diagnóstico.1 2 O delírio é comum no ambiente de cuidados agudos, incluindo o departamento de emergência (DE), e é particularmente prev- alente entre adultos com mais de 65 anos de idade. Até 7%-17% dos adultos mais velhos que se apresentam à DE preenchem os critérios de diagnóstico para o delírio.3-9 Os prestadores de "
$new = "O delírio é um estado de confusão aguda caracterizado por declínios na atenção, na consciência e na cognição. As flutuações do estado mental ao longo do tempo são características e necessárias para o diagnóstico.1 2 O delírio é comum no ambiente de cuidados agudos, incluindo o departamento de emergência (DE), e é particularmente prevalente entre adultos com mais de 65 anos de idade. Até 7%-17% dos adultos mais velhos que se apresentam à DE preenchem os critérios de diagnóstico para o delírio.3-9 Os prestadores de "

$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
